# CharacterSkill.xlsx — "Change Structure & Move Cur Unit to Level Data"
#
# Restructures the skill table: adds "range"/"radius" (int) columns,
# renumbers skill ids from short codes to level-scoped codes, renames the
# first skill and gives the second skill a name, and drops the old,
# now-unused third data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Best-effort: match the author's resized window (cosmetic only) ---
$excel.Width = 25600
$excel.Height = 12210

# --- New columns D:E — header row + "int" type row ---
$ws.Range("D2").Value = "range"
$ws.Range("E2").Value = "radius"
$ws.Range("D3").Value = "int"
$ws.Range("E3").Value = "int"

# --- Row 4: skill 1001 "Water" -> 11001 "Elemental Bolt" w/ range 4, radius 1 ---
$ws.Range("A4").Value = 11001
$ws.Range("B4").Value = "Elemental Bolt"
$ws.Range("C4").ClearContents()
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 1

# --- Row 5: skill 2001 -> 21001 "Strike" w/ range 1, radius 1 ---
$ws.Range("A5").Value = 21001
$ws.Range("B5").Value = "Strike"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# --- Old row 6 (skill 3001, unused) is removed; table now ends at row 5 ---
$ws.Rows.Item(6).Delete()

# --- Leave selection where the author's session ended ---
$ws.Range("E5").Select()
